# Update the TableComparator summary numbers on Sheet1.
#
# The underlying comparison logic was fixed (generic BigDecimal handling),
# which reclassifies one id from "NO MATCH" into "MATCH":
#   - Quality global numerador/pct:  0 / 0.0%   -> 1 / 10.0%
#   - MATCH          numerador/pct/ejemplos: 1 / 14.3% / NULL        -> 2 / 28.6% / 1,NULL
#   - NO MATCH       numerador/pct/ejemplos: 6 / 85.7% / 1,2,4,7,8,9 -> 5 / 71.4% / 2,4,7,8,9
#
# All of these cells store plain text (not real numbers/percentages) in the
# original workbook, so we force Excel to keep them as text (NumberFormat
# "@") instead of letting it auto-convert "1" / "10.0%" into numeric values,
# then strip the temporary formatting back off so no stray cell style is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 7: "Quality global"
Set-TextValue $ws.Range("D7") "1"
Set-TextValue $ws.Range("F7") "10.0%"

# Row 8: "MATCH" / "1:1 (exact matches)"
Set-TextValue $ws.Range("D8") "2"
Set-TextValue $ws.Range("F8") "28.6%"
$ws.Range("G8").Value = "1,NULL"

# Row 9: "NO MATCH" / "1:1 (match not identical)"
Set-TextValue $ws.Range("D9") "5"
Set-TextValue $ws.Range("F9") "71.4%"
$ws.Range("G9").Value = "2,4,7,8,9"
